# Apply the "add RMSE percentage for Morocco soil" edit:
#  - Rename header B1 "RMSE" -> "Yield RMSE"
#  - Add new headers C1, D1, E1
#  - Add new data columns C (Yield RMSE Percentage), D (Water Used RMSE),
#    E (Water Used RMSE Percentage) for rows 2-15
#  - Tiny floating point refresh on B6, B9, B12 (last digit changed)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("B1").Value = "Yield RMSE"
$ws.Range("C1").Value = "Yield RMSE Percentage"
$ws.Range("D1").Value = "Water Used RMSE"
$ws.Range("E1").Value = "Water Used RMSE Percentage"

# Copy the header style (bold/border/center/top) from B1 into the new header cells
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1:E1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Re-apply values since PasteSpecial(formats) shouldn't touch them, but make sure
$ws.Range("C1").Value = "Yield RMSE Percentage"
$ws.Range("D1").Value = "Water Used RMSE"
$ws.Range("E1").Value = "Water Used RMSE Percentage"

# --- Data rows ---
# Correct tiny floating point differences in column B
$ws.Range("B6").Value = 1.536866645809522
$ws.Range("B9").Value = 1.261615158411222
$ws.Range("B12").Value = 1.775887473276425

# Column C: Yield RMSE Percentage
$yieldPct = @{
    2  = 50.66685706577445
    3  = 68.73412415407923
    4  = 68.88435675727116
    5  = 68.61642608117414
    6  = 71.87807081939918
    7  = 63.22341356323691
    8  = 64.18244098157632
    9  = 59.00477048569406
    10 = 103.0147907626998
    11 = 77.05674707265599
    12 = 83.05689105785127
    13 = 73.1865429698023
    14 = 88.68898557385536
    15 = 87.36650856204157
}

# Column D: Water Used RMSE
$waterRmse = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 10.88214375165017
}

# Column E: Water Used RMSE Percentage
$waterPct = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 6.691285801985544
}

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 3).Value = $yieldPct[$r]
    $ws.Cells.Item($r, 4).Value = $waterRmse[$r]
    $ws.Cells.Item($r, 5).Value = $waterPct[$r]
}
